# Edit description (from commit diff):
#   1. Three tables (on slides 14, 15 and 16) switch from the custom
#      "Table_0" table style to the built-in table style
#      {BA71F996-A943-475B-B67B-8CABE722F230}.
#   2. The deck's theme colour palette changes from the "Integral" /
#      "Red Violet" scheme to the standard Office "Office" colour scheme
#      (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# Table styles can't be changed by assigning Table.Style directly (that
# throws - PowerPoint requires Table.ApplyStyle(id) instead), so every
# table on the three affected slides is located and re-styled that way.
#
# The colour palette is edited through the slide's ThemeColorScheme,
# which is the live view onto the presentation's shared theme - setting
# each of the twelve colour slots updates the theme for the whole deck.

$p = $ppt.ActivePresentation

$newTableStyleId = "{BA71F996-A943-475B-B67B-8CABE722F230}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# Office theme colour scheme, in ThemeColorScheme slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($k = 1; $k -le $officeColors.Length; $k++) {
    $tcs.Item($k).RGB = $officeColors[$k - 1]
}
